$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (capital D) to match the workbook's sheet tag rename
$ws.Name = "Dados.xlsx"

# Fix header typo: "Nota Ciencia" -> "Nota Ciências"
$ws.Range("E1").Value = "Nota Ciências"

# Insert a new row 2 for "Claudio", pushing João/Maria/Pedro/Ana down by one
$ws.Rows.Item(2).Insert()

# Make sure all the data cells (B:E for every student row) are stored as
# text (not numbers), matching the target workbook's inlineStr cells.
$ws.Range("B2:E6").NumberFormat = "@"

# Row 2: Claudio (new row)
$ws.Range("A2").Value = "Claudio"
$ws.Range("B2").Value = "28"
$ws.Range("C2").Value = "8"
$ws.Range("D2").Value = "7"
$ws.Range("E2").Value = "9"
$ws.Range("F2").Value = "Insuficiente"
$ws.Range("G2").Value = "Insuficiente"
$ws.Range("H2").Value = "Insuficiente"

# Row 3: João (values unchanged from before, but re-typed as text)
$ws.Range("A3").Value = "João"
$ws.Range("B3").Value = "16"
$ws.Range("C3").Value = "8"
$ws.Range("D3").Value = "7"
$ws.Range("E3").Value = "6"
$ws.Range("F3").Value = "Insuficiente"
$ws.Range("G3").Value = "Insuficiente"
$ws.Range("H3").Value = "Insuficiente"

# Row 4: Maria (C/D values changed)
$ws.Range("A4").Value = "Maria"
$ws.Range("B4").Value = "17"
$ws.Range("C4").Value = "19"
$ws.Range("D4").Value = "12"
$ws.Range("E4").Value = "14"
$ws.Range("F4").Value = "Suficiente"
$ws.Range("G4").Value = "Suficiente"
$ws.Range("H4").Value = "Suficiente"

# Row 5: Pedro (values unchanged from before, but re-typed as text)
$ws.Range("A5").Value = "Pedro"
$ws.Range("B5").Value = "16"
$ws.Range("C5").Value = "6"
$ws.Range("D5").Value = "9"
$ws.Range("E5").Value = "10"
$ws.Range("F5").Value = "Insuficiente"
$ws.Range("G5").Value = "Insuficiente"
$ws.Range("H5").Value = "Suficiente"

# Row 6: Ana (values unchanged from before, but re-typed as text)
$ws.Range("A6").Value = "Ana"
$ws.Range("B6").Value = "18"
$ws.Range("C6").Value = "15"
$ws.Range("D6").Value = "20"
$ws.Range("E6").Value = "17"
$ws.Range("F6").Value = "Suficiente"
$ws.Range("G6").Value = "Suficiente"
$ws.Range("H6").Value = "Suficiente"
